$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-07-06 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-07 Sunday", 2)

$d.Content.Find.Execute("834÷9=92, 6", $true, $false, $false, $false, $false, $true, 1, $false, "238÷6=39, 4", 2)
$d.Content.Find.Execute("821÷5=164, 1", $true, $false, $false, $false, $false, $true, 1, $false, "318÷9=35, 3", 2)
$d.Content.Find.Execute("339÷6=56, 3", $true, $false, $false, $false, $false, $true, 1, $false, "816÷6=136, 0", 2)
$d.Content.Find.Execute("215÷6=35, 5", $true, $false, $false, $false, $false, $true, 1, $false, "920÷2=460, 0", 2)
$d.Content.Find.Execute("368÷6=61, 2", $true, $false, $false, $false, $false, $true, 1, $false, "858÷8=107, 2", 2)

$d.Content.Find.Execute("145÷8=18, 1", $true, $false, $false, $false, $false, $true, 1, $false, "997÷5=199, 2", 2)
$d.Content.Find.Execute("748÷9=83, 1", $true, $false, $false, $false, $false, $true, 1, $false, "422÷8=52, 6", 2)
$d.Content.Find.Execute("970÷2=485, 0", $true, $false, $false, $false, $false, $true, 1, $false, "154÷9=17, 1", 2)
$d.Content.Find.Execute("555÷6=92, 3", $true, $false, $false, $false, $false, $true, 1, $false, "613÷5=122, 3", 2)
$d.Content.Find.Execute("473÷6=78, 5", $true, $false, $false, $false, $false, $true, 1, $false, "803÷9=89, 2", 2)

$d.Content.Find.Execute("634÷4=158, 2", $true, $false, $false, $false, $false, $true, 1, $false, "734÷4=183, 2", 2)
$d.Content.Find.Execute("478÷8=59, 6", $true, $false, $false, $false, $false, $true, 1, $false, "174÷6=29, 0", 2)
$d.Content.Find.Execute("758÷6=126, 2", $true, $false, $false, $false, $false, $true, 1, $false, "811÷7=115, 6", 2)
$d.Content.Find.Execute("524÷4=131, 0", $true, $false, $false, $false, $false, $true, 1, $false, "452÷5=90, 2", 2)
$d.Content.Find.Execute("680÷5=136, 0", $true, $false, $false, $false, $false, $true, 1, $false, "567÷6=94, 3", 2)

$d.Content.Find.Execute("868÷3=289, 1", $true, $false, $false, $false, $false, $true, 1, $false, "711÷4=177, 3", 2)
$d.Content.Find.Execute("599÷3=199, 2", $true, $false, $false, $false, $false, $true, 1, $false, "195÷4=48, 3", 2)
$d.Content.Find.Execute("240÷7=34, 2", $true, $false, $false, $false, $false, $true, 1, $false, "146÷9=16, 2", 2)
$d.Content.Find.Execute("424÷7=60, 4", $true, $false, $false, $false, $false, $true, 1, $false, "802÷2=401, 0", 2)
$d.Content.Find.Execute("105÷8=13, 1", $true, $false, $false, $false, $false, $true, 1, $false, "459÷9=51, 0", 2)

$d.Content.Find.Execute("388÷8=48, 4", $true, $false, $false, $false, $false, $true, 1, $false, "963÷9=107, 0", 2)
$d.Content.Find.Execute("139÷8=17, 3", $true, $false, $false, $false, $false, $true, 1, $false, "780÷6=130, 0", 2)
$d.Content.Find.Execute("165÷6=27, 3", $true, $false, $false, $false, $false, $true, 1, $false, "578÷4=144, 2", 2)
$d.Content.Find.Execute("822÷8=102, 6", $true, $false, $false, $false, $false, $true, 1, $false, "424÷9=47, 1", 2)
$d.Content.Find.Execute("276÷2=138, 0", $true, $false, $false, $false, $false, $true, 1, $false, "302÷2=151, 0", 2)
